$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10: moonshotai/kimi-k2-instruct-0905
$ws.Range("I10").Value = 2
$ws.Range("J10").Value = 0.002
$ws.Range("K10").Value = 945
$ws.Range("L10").Value = 0.00315

# Row 14: qwen/qwen3-32b
$ws.Range("I14").Value = 2
$ws.Range("J14").Value = 0.002
$ws.Range("K14").Value = 939
$ws.Range("L14").Value = 0.001878
